$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simple scalar values on row 5
$ws.Range("B5").Value = 80221
$ws.Range("E5").Value = 389
$ws.Range("F5").Value = "Läderlappslav"
$ws.Range("G5").Value = "Collema nigrescens"
$ws.Range("H5").Value = "(Huds.) DC."

# Create new (empty) cells at J5, K5, N5 to mirror the blank inline-string
# placeholder cells that already exist elsewhere in the sheet (e.g. I5).
# A harmless formatting touch forces Excel to materialize a cell entry
# without giving it any actual value/content.
$ws.Range("J5").Font.Bold = $false
$ws.Range("K5").Font.Bold = $false
$ws.Range("N5").Font.Bold = $false

# Flip the boolean flag in AE5
$ws.Range("AE5").Value = $true

# Create a new empty cell at AF5 as well
$ws.Range("AF5").Font.Bold = $false
